$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "'9"
$ws.Range("D3").Value = "'23.06"
$ws.Range("G3").Value = "'9"
$ws.Range("D4").Value = "'5.409"
$ws.Range("G4").Value = "'9"
$ws.Range("D5").Value = "'0.05920"
$ws.Range("G5").Value = "'9"
$ws.Range("G6").Value = "'9"
$ws.Range("D7").Value = "'6.551"
$ws.Range("G7").Value = "'9"
$ws.Range("D8").Value = "'0.8128"
$ws.Range("G8").Value = "'9"
$ws.Range("D9").Value = "'0.9111"
$ws.Range("G9").Value = "'9"
$ws.Range("D10").Value = "'0.1404"
$ws.Range("G10").Value = "'9"
$ws.Range("D11").Value = "'0.07407"
$ws.Range("G11").Value = "'9"
$ws.Range("D12").Value = "'0.03275"
$ws.Range("G12").Value = "'9"
$ws.Range("D13").Value = "'0.03052"
$ws.Range("G13").Value = "'9"
$ws.Range("D14").Value = "'0.09351"
$ws.Range("G14").Value = "'9"
$ws.Range("D15").Value = "'3.851"
$ws.Range("G15").Value = "'9"
$ws.Range("D16").Value = "'0.001579"
$ws.Range("G16").Value = "'9"
$ws.Range("D17").Value = "'0.04675"
$ws.Range("G17").Value = "'9"
$ws.Range("G18").Value = "'9"
$ws.Range("D19").Value = "'0.006114"
$ws.Range("G19").Value = "'9"
$ws.Range("D20").Value = "'0.004977"
$ws.Range("G20").Value = "'9"
$ws.Range("D21").Value = "'0.0009853"
$ws.Range("G21").Value = "'9"
$ws.Range("D22").Value = "'0.00009406"
$ws.Range("G22").Value = "'9"
$ws.Range("D23").Value = "'3.604"
$ws.Range("G23").Value = "'9"
$ws.Range("G24").Value = "'9"
$ws.Range("D25").Value = "'0.3238"
$ws.Range("G25").Value = "'9"
$ws.Range("G26").Value = "'9"
$ws.Range("D27").Value = "'0.0001503"
$ws.Range("G27").Value = "'9"
$ws.Range("G28").Value = "'9"
$ws.Range("G29").Value = "'9"
$ws.Range("G30").Value = "'9"
$ws.Range("G31").Value = "'9"
$ws.Range("G32").Value = "'9"
$ws.Range("G33").Value = "'9"
$ws.Range("G34").Value = "'9"
$ws.Range("G35").Value = "'9"
$ws.Range("G36").Value = "'9"
$ws.Range("G37").Value = "'9"
$ws.Range("G38").Value = "'9"
$ws.Range("G39").Value = "'9"
$ws.Range("D40").Value = "'0.03990"
$ws.Range("G40").Value = "'9"
$ws.Range("D41").Value = "'0.006203"
$ws.Range("G41").Value = "'9"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("G42").Value = "'9"
$ws.Range("D43").Value = "'0.003002"
$ws.Range("G43").Value = "'9"
$ws.Range("D44").Value = "'0.008860"
$ws.Range("G44").Value = "'9"
$ws.Range("D45").Value = "'0.00005248"
$ws.Range("G45").Value = "'9"
$ws.Range("G46").Value = "'9"
$ws.Range("B47").Value = "'CoinbaseStockToken"
$ws.Range("C47").Value = "'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "'0.7823"
$ws.Range("E47").Value = "'46CoinbaseStockTokenCOIN"
$ws.Range("G47").Value = "'9"
$ws.Range("B48").Value = "'BOLO"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.002264"
$ws.Range("E48").Value = "'47BOLOBOLO"
$ws.Range("G48").Value = "'9"
$ws.Range("B49").Value = "'CryptobidCoin"
$ws.Range("C49").Value = "'https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'48CryptobidCoinCBC"
$ws.Range("G49").Value = "'9"
$ws.Range("B50").Value = "'SpecialPowerGold"
$ws.Range("C50").Value = "'https://coinranking.com/coin/jPTWzmsWb+specialpowergold-spg"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'49SpecialPowerGoldSPGWorstin24h"
$ws.Range("G50").Value = "'9"
$ws.Range("B51").Value = "'DigiFinexToken"
$ws.Range("C51").Value = "'https://coinranking.com/coin/rY6dWXQL4+digifinextoken-dft"
$ws.Range("D51").Value = "'--"
$ws.Range("E51").Value = "'50DigiFinexTokenDFT"
$ws.Range("G51").Value = "'9"
